# Applies the row-2 data update on the "IssueAPIData" sheet of TestData.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("IssueAPIData")

$ws.Range("B2").Value = "create new issue "
$ws.Range("C2").Value = "create new issue description"
$ws.Range("D2").Value = 1234
$ws.Range("E2").Value = "regression"
$ws.Range("F2").Value = "issue"
$ws.Range("G2").Value = "open"
